# v0.2: Square can have a background.
$wb = $excel.ActiveWorkbook

$wsMaze = $wb.Worksheets.Item(1)
$wsBg   = $wb.Worksheets.Item(2)

# --- Rename sheets ---------------------------------------------------------
$wsMaze.Name = "maze"
$wsBg.Name   = "background"

# --- background sheet: mark squares that carry a background with "bg" -----
# Existing single-letter/legend cells are replaced by a uniform "bg" marker,
# re-using B2's cell style (centered, unbordered) for every marker cell.
$wsBg.Range("C3").Clear()
$wsBg.Range("D4").Clear()
$wsBg.Range("B5").Clear()
$wsBg.Range("E6").Clear()

$wsBg.Range("B2").Value = "bg"
$wsBg.Range("C4").Value = "bg"

$wsBg.Range("B2").Copy() | Out-Null
$wsBg.Range("D5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsBg.Range("D5").Value = "bg"

$wsBg.Range("D7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsBg.Range("D7").Value = "bg"
$excel.CutCopyMode = $false

# --- match the maze sheet's row height on the background sheet ------------
for ($r = 2; $r -le 8; $r++) {
  $wsBg.Rows.Item($r).RowHeight = 39.55
}

# --- selections -------------------------------------------------------------
$wsMaze.Range("G11").Select()
$wsBg.Activate()
$wsBg.Range("D8").Select()
